$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - logistic_embeddings
$ws.Range("C5").Value = 0.624
$ws.Range("D5").Value = 0.727
$ws.Range("E5").Value = 0.753
$ws.Range("F5").Value = 0.789
$ws.Range("G5").Value = 0.681
$ws.Range("H5").Value = 0.695

# Row 7 - classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.624
$ws.Range("E7").Value = 0.753
$ws.Range("F7").Value = 0.789
$ws.Range("G7").Value = 0.681
$ws.Range("H7").Value = 0.695

# Row 8 - BERT-base
$ws.Range("C8").Value = 0.654
$ws.Range("D8").Value = 0.761
$ws.Range("E8").Value = 0.788
$ws.Range("F8").Value = 0.794
$ws.Range("G8").Value = 0.727
$ws.Range("H8").Value = 0.747

# Row 9 - BERT-base-nli
$ws.Range("B9").Value = 0.555
$ws.Range("C9").Value = 0.628
$ws.Range("D9").Value = 0.745
$ws.Range("E9").Value = 0.77
$ws.Range("F9").Value = 0.771
$ws.Range("G9").Value = 0.678
$ws.Range("H9").Value = 0.6919999999999999
